$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference a style from a row untouched by this update (row 45) to strip the
# quote-prefix style Excel applies when a number-looking value is forced to text,
# keeping the re-written cells on the same (unstyled) style as the rest of column D.
$plainStyle = $ws.Range("D45").Style

$ws.Range("D2").Value = "91.630.10"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").Value = "3.083.77"
$ws.Range("E3").Value = "  -2.51%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'234.79"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  -2.39%  "

$ws.Range("D6").Value = "'608.93"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  -1.93%  "

$ws.Range("E7").Value = "  -4.23%  "

$ws.Range("D8").Value = "'0.382"
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = "  +1.32%  "

$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").Value = "3.079.77"
$ws.Range("E10").Value = "  -2.53%  "

$ws.Range("E11").Value = "  -3.29%  "

$ws.Range("E12").Value = "  -1.88%  "

$ws.Range("D13").Value = "'0.0000245"
$ws.Range("D13").Style = $plainStyle
$ws.Range("E13").Value = "  -1.56%  "

$ws.Range("D14").Value = "91.988.17"
$ws.Range("E14").Value = "  +0.80%  "

$ws.Range("D15").Value = "'33.73"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "  -4.65%  "

$ws.Range("D16").Value = "'5.38"
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = "  -4.01%  "

$ws.Range("D17").Value = "3.668.95"
$ws.Range("E17").Value = "  -2.06%  "

$ws.Range("D18").Value = "3.094.07"
$ws.Range("E18").Value = "  -2.56%  "

$ws.Range("D19").Value = "'3.65"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = "  -3.20%  "

$ws.Range("D20").Value = "'14.52"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "  -4.85%  "

$ws.Range("D21").Value = "'5.72"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  -5.96%  "

$ws.Range("D22").Value = "'9.22"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").Value = "'440.74"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  -3.68%  "

$ws.Range("D24").Value = "'0.0000192"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  -6.40%  "

$ws.Range("D25").Value = "'5.67"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = "  -5.91%  "

$ws.Range("D26").Value = "'85.57"
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = "  -4.15%  "

$ws.Range("D27").Value = "'11.51"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  -4.90%  "

$ws.Range("D28").Value = "3.257.67"
$ws.Range("E28").Value = "  -2.02%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").Value = "'0.129"
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = "  -2.26%  "

$ws.Range("D31").Value = "'0.228"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = "  -3.13%  "

$ws.Range("E32").Value = "  -2.78%  "

$ws.Range("D33").Value = "'9.03"
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = "  -4.33%  "

$ws.Range("D34").Value = "'0.993"
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = "  -0.76%  "

$ws.Range("D35").Value = "'7.75"
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = "  +0.25%  "

$ws.Range("D36").Value = "'0.157"
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = "  -9.00%  "

$ws.Range("D37").Value = "'25.73"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = "  -3.20%  "

$ws.Range("D38").Value = "'3.89"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "  -0.82%  "

$ws.Range("D39").Value = "'1.88"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  -4.48%  "

$ws.Range("D40").Value = "'479.51"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = "  -6.92%  "

$ws.Range("D41").Value = "'23.86"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  +7.74%  "

$ws.Range("E42").Value = "  -6.16%  "

$ws.Range("D43").Value = "'0.428"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  -5.58%  "

$ws.Range("D44").Value = "'3.26"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  -6.30%  "

$ws.Range("D46").Value = "'162.81"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  +3.21%  "

# Rows 47/48 swap places: ARBITRUM <-> Stacks (with refreshed price/volume).
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'1.86"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  -5.06%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'0.680"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  -5.45%  "

$ws.Range("E49").Value = "  -1.23%  "

$ws.Range("D50").Value = "'0.0330"
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = "  +1.43%  "

$ws.Range("D51").Value = "'43.83"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  -0.48%  "
